$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("A2").Value = 98
$ws.Range("B2").Value = 194
$ws.Range("C2").Value = 137

# Row 3 updates (previously empty cells)
$ws.Range("A3").Value = 19
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 41
$ws.Range("J3").Formula = "=13+16+31"

# Row 4 updates (previously empty cells)
$ws.Range("A4").Value = 143
$ws.Range("B4").Value = 300
$ws.Range("C4").Value = 205
$ws.Range("F4").Value = 55
$ws.Range("G4").Value = 115
$ws.Range("H4").Value = 83

# Row 5 updates (previously empty cells)
$ws.Range("A5").Value = 81
$ws.Range("B5").Value = 150
$ws.Range("C5").Value = 91

# Row 6 updates (previously empty cells)
$ws.Range("A6").Value = 31
$ws.Range("B6").Value = 60
$ws.Range("C6").Value = 40

# Update selection to match target
$ws.Range("C7").Select()
